$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update label for row 9 (August cutoff date changed)
$ws.Range("A9").Value = "August (through 08-26)"

# Update row 9 values (August)
$ws.Range("B9").Value = 27
$ws.Range("C9").Value = 63
$ws.Range("D9").Value = 75
$ws.Range("E9").Value = 52
$ws.Range("F9").Value = 38
$ws.Range("G9").Value = 145
$ws.Range("H9").Value = 132

# Update row 10 values (Total)
$ws.Range("B10").Value = 189
$ws.Range("C10").Value = 365
$ws.Range("D10").Value = 540
$ws.Range("E10").Value = 477
$ws.Range("F10").Value = 342
$ws.Range("G10").Value = 766
$ws.Range("H10").Value = 1046
